$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add the new "VendorMaster" sheet, positioned right before
#    "IndirectRequisition" (i.e. right after "NONINVPOITEM").
# ------------------------------------------------------------------
$indirectReq = $wb.Worksheets.Item("IndirectRequisition")
$vendorMaster = $wb.Worksheets.Add($indirectReq)
$vendorMaster.Name = "VendorMaster"

# Fill in the header row / data row in the same order the strings were
# introduced into the shared-string table (Street..Country, then the
# second data row, and finally the Buyer column last).
$vendorMaster.Range("B1").Value = "Street"
$vendorMaster.Range("C1").Value = "City"
$vendorMaster.Range("D1").Value = "State"
$vendorMaster.Range("E1").Value = "Zip"
$vendorMaster.Range("F1").Value = "Country"

$vendorMaster.Range("B2").Value = "Avenue Villa"
$vendorMaster.Range("C2").Value = "Boston"
$vendorMaster.Range("D2").Value = "Massachusetts"
$vendorMaster.Range("F2").Value = "US"

$vendorMaster.Range("A1").Value = "Buyer"
$vendorMaster.Range("A2").Value = "Namrata Patil"
$vendorMaster.Range("E2").Value = 2101

$vendorMaster.Columns.Item(1).AutoFit() | Out-Null

# ------------------------------------------------------------------
# 2. POReceipt: remove the stray "il..." note in G10 (it becomes an
#    orphan shared string that drops out of the table on save) and
#    make POReceipt the active sheet/selection.
# ------------------------------------------------------------------
$poReceipt = $wb.Worksheets.Item("POReceipt")
$poReceipt.Range("G10").ClearContents()

# ------------------------------------------------------------------
# 3. Update the active sheet / selections to match the saved view
#    state: POReceipt becomes the active tab with G15 selected, and
#    the new VendorMaster sheet is left with E1 selected (not active).
# ------------------------------------------------------------------
$vendorMaster.Activate()
$vendorMaster.Range("E1").Select()

$poReceipt.Activate()
$poReceipt.Range("G15").Select()
